$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC40_Adding_MultipleItems_Quick")

# Insert a new row at position 93, shifting rows 93:95 down to 94:96
$ws.Rows.Item(93).Insert()

# New row 93 content: B93 = WAIT (matches string used in row 92/90/88)
$ws.Cells.Item(93, 2).Value = "WAIT"

# Update styles on B91 and B92 to style index 2 (fill + border, no font)
$ws.Cells.Item(91, 2).Style = "Normal"
$ws.Range("B91").Interior.ColorIndex = 13
$ws.Range("B91").Borders.LineStyle = 1
$ws.Range("B92").Interior.ColorIndex = 13
$ws.Range("B92").Borders.LineStyle = 1

# Selection
$ws.Range("B90:B94").Select()

$wb.Save()
